$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions): bump the "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 201
$ws1.Range("F3").Value = 5500
$ws1.Range("F12").Value = 4971
$ws1.Range("F21").Value = 115

# Sheet "全部类型" (all types): same rows mirrored, different row offsets
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 201
$ws4.Range("F4").Value = 5500
$ws4.Range("F13").Value = 4971
$ws4.Range("F22").Value = 115
